# Apply the "in-out change" edit to data_example.xlsx
#
# Summary of changes (from the authoritative diff):
#  - device sheet: add new column P "p_ratio" with header + values
#       P1 = "p_ratio", P4 = 1, P5 = 1, P6 = 1
#  - device sheet: I6 0.95 -> 0.9
#  - device sheet: N7 10 -> 100, O7 10 -> 100
#  - device sheet: O8 300 -> 200
#  - device sheet: O10 10 -> 0
#  - device sheet: O12 250 -> 100
#  - device sheet: O13 10 -> 0
#  - device sheet: O14 0 -> 10
#  - device sheet: C15 0 -> 1
#  - selection on "edge" sheet moved to C9
#  - selection on "device" sheet moved to O8 (and it's the active tab)

$wb = $excel.ActiveWorkbook

$wsEdge   = $wb.Worksheets.Item("edge")
$wsDevice = $wb.Worksheets.Item("device")

# --- device sheet: new "p_ratio" column (P) ---
$wsDevice.Range("P1").Value = "p_ratio"
$wsDevice.Range("P4").Value = 1
$wsDevice.Range("P5").Value = 1
$wsDevice.Range("P6").Value = 1

# --- device sheet: value updates ---
$wsDevice.Range("I6").Value = 0.9

$wsDevice.Range("N7").Value = 100
$wsDevice.Range("O7").Value = 100

$wsDevice.Range("O8").Value = 200

$wsDevice.Range("O10").Value = 0

$wsDevice.Range("O12").Value = 100

$wsDevice.Range("O13").Value = 0

$wsDevice.Range("O14").Value = 10

$wsDevice.Range("C15").Value = 1

# --- selections ---
$wsEdge.Select()
$wsEdge.Range("C9").Select()

$wsDevice.Select()
$wsDevice.Range("O8").Select()
